# Applies the "Hydrogen Energy" -> "Chemistry" rewrite described by the
# supplied diff: title, byline/author, e-mail, all body paragraphs, the
# summary paragraph, and a trailing empty paragraph appended at the end.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find/Replace failed for: $old"
    }
}

# --- Title ---------------------------------------------------------------
Replace-Text "Hydrogen Energy: A Sustainable Frontier" "The Marvelous World of Chemistry: Uncovering the Secrets of Matter"

# --- Byline (Emily Watson -> Dr. Amelia Carter) ---------------------------
Replace-Text "Emily Watson" "Dr. Amelia Carter"

# --- E-mail address --------------------------------------------------------
# Original runs: "emily" | "." | "watson@emailworld" | "." | "com"
# Target runs:    "ameliacartere@schoolmail" | "." | "edu"
Replace-Text "emily" "ameliacartere@schoolmail"
Replace-Text "watson@emailworld.com" "edu"

# --- Body paragraph 1 (intro) ---------------------------------------------
Replace-Text "The quest for clean, renewable energy sources has intensified as the world faces the twin challenges of climate change and diminishing fossil fuels" "The realm of chemistry is filled with countless wonders and mysteries waiting to be unraveled"
Replace-Text " Among the promising contenders, hydrogen energy stands out with its immense potential to transform the global energy landscape" " It is a discipline that pushes the boundaries of human knowledge, constantly revealing new insights into the nature of matter"
Replace-Text " This versatile fuel offers a pathway toward sustainable development, fostering energy security, reducing carbon emissions, and revolutionizing industries across the spectrum" " The study of chemistry not only enhances our understanding of the world but also equips us with the skills and knowledge to solve real-world problems"

# --- Body paragraph 2 (properties) -----------------------------------------
Replace-Text "Hydrogen, the lightest and most abundant element, possesses a remarkable energy density" "Chemistry is an integral part of our educational system, providing a solid foundation for further studies in STEM (Science, Technology, Engineering, and Mathematics) fields"
Replace-Text " When combined with oxygen, it undergoes combustion, releasing substantial heat and producing water as a byproduct, thereby achieving zero carbon emissions. This clean-burning fuel holds the key to decarbonizing sectors like transportation, heating, and electricity generation, which collectively contribute significantly to greenhouse gas emissions" " It prepares students for careers in various industries, ranging from pharmaceuticals and biotechnology to materials science and environmental engineering"

# --- Body paragraph 3 (production) -----------------------------------------
Replace-Text "The production of hydrogen, however, presents a challenge" "Whether you aspire to become a chemist, a doctor, an engineer, or a teacher, a strong foundation in chemistry will serve as an invaluable asset"
Replace-Text " Traditional methods rely on fossil fuels, negating the environmental benefits" " It will open up a world of opportunities and empower you to make a positive impact on society"
Replace-Text " The focus is now shifting toward green hydrogen, produced from renewable sources like solar and wind energy. This zero-carbon process, known as electrolysis, splits water molecules into hydrogen and oxygen, paving the way for a sustainable hydrogen economy" " So, embrace the wonders of chemistry, embark on a journey of discovery, and unravel the secrets of the universe"

# --- Summary paragraph -------------------------------------------------------
Replace-Text "Hydrogen energy has emerged as a promising solution to the global energy crisis" "Chemistry is a captivating field that delves into the structure, properties, and interactions of matter"
Replace-Text " Its versatility, clean-burning nature, and zero-carbon emissions make it a strong contender for decarbonizing industries and achieving sustainable development" " It plays a pivotal role in shaping our understanding of the universe, unlocking the mysteries of life, and addressing global challenges"
Replace-Text " While the production of hydrogen remains a challenge, advancements in green hydrogen generation from renewable sources are driving progress toward a hydrogen economy" " Chemistry offers a window into the microscopic world, enabling us to manipulate matter at the atomic and molecular levels"
Replace-Text " As technology evolves, hydrogen energy holds immense potential to transform the way we power our world and secure a sustainable future" " Its applications are vast, ranging from the development of new materials and medicines to the advancement of renewable energy sources. The study of chemistry empowers us with the knowledge and skills to address complex problems, unlocking a world of opportunities and enabling us to make a positive impact on society"

# --- Trailing empty paragraph ------------------------------------------------
$last = $d.Paragraphs($d.Paragraphs.Count)
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()

Write-Host "All replacements applied."
